$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '41.656.82'
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.469.77'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('E4').Value = '  +0.13%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '316.74'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.53%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '92.26'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -1.39%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.553'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +1.91%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.515'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +1.77%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0892'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +13.24%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '32.68'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('E12').Value = '  -0.51%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '2.848.38'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -0.71%  '
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('E15').Value = '  -2.66%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '2.471.69'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.58%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.786'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +3.29%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '41.600.41'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +0.11%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.0₃0963'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +2.41%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '6.46'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +1.33%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '71.43'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.20%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '11.53'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.81%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '241.01'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +1.52%  '
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').Value = '  +0.00%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '24.83'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('E28').Value = '  +3.59%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.84'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.36%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '35.36'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -1.99%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '156.18'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -0.03%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.0767'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +1.82%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '2.57'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.06%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '17.45'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -1.71%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.90'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -1.35%  '
$ws.Range('E38').Value = '  +1.35%  '
$ws.Range('E39').Value = '  -2.20%  '
$ws.Range('E40').Value = '  -2.48%  '
$ws.Range('B41').Value = 'ApeXProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.50'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.29%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '3.96'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -3.39%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.979.92'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '19.00'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -4.77%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0284'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.46%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.99'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('E47').Value = '  +1.37%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.704.24'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.70%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '97.15'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -0.51%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '74.23'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.95%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '67.09'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -2.22%  '
